# Auto-generated: apply Kraken_Profits-style profit recompute to each leve sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 506
$ws.Range("I2").Value = 506
$ws.Range("K2").Value = 506
$ws.Range("M2").Value = -393
# Row 4
$ws.Range("H4").Value = 3896.8462
$ws.Range("I4").Value = 2532.6365
$ws.Range("K4").Value = 2532.6365
$ws.Range("M4").Value = -2418.6365
# Row 5
$ws.Range("H5").Value = 220.9
$ws.Range("I5").Value = 189.88889
$ws.Range("K5").Value = 189.88889
$ws.Range("M5").Value = -74.88889
# Row 43
$ws.Range("H43").Value = 3299
$ws.Range("J43").Value = 2999.5
$ws.Range("L43").Value = 2999.5
$ws.Range("N43").Value = -3137.5
# Row 70
$ws.Range("H70").Value = 27599.75
$ws.Range("J70").Value = 36333
$ws.Range("L70").Value = 108999
$ws.Range("N70").Value = -109539
# Row 73
$ws.Range("H73").Value = 27599.75
$ws.Range("J73").Value = 36333
$ws.Range("L73").Value = 108999
$ws.Range("N73").Value = -110871
# Row 98
$ws.Range("H98").Value = 1766.5
$ws.Range("I98").Value = 1485.2307
$ws.Range("J98").Value = 2985.3333
$ws.Range("K98").Value = 1485.2307
$ws.Range("L98").Value = 2985.3333
$ws.Range("M98").Value = 12.76929999999993
$ws.Range("N98").Value = -5981.3333
# Row 111
$ws.Range("H111").Value = 209.66667
$ws.Range("I111").Value = 209.66667
$ws.Range("K111").Value = 629.00001
$ws.Range("M111").Value = 2437.99999
# Row 118
$ws.Range("H118").Value = 4089
$ws.Range("I118").Value = 445
$ws.Range("K118").Value = 1335
$ws.Range("M118").Value = 322
# Row 122
$ws.Range("H122").Value = 1766.5
$ws.Range("I122").Value = 1485.2307
$ws.Range("J122").Value = 2985.3333
$ws.Range("K122").Value = 4455.6921
$ws.Range("L122").Value = 8955.999899999999
$ws.Range("M122").Value = -2005.6921
$ws.Range("N122").Value = -13855.9999
# Row 135
$ws.Range("H135").Value = 1368.909
$ws.Range("I135").Value = 1139.7778
$ws.Range("K135").Value = 10258.0002
$ws.Range("M135").Value = -7723.0002
# Row 137
$ws.Range("H137").Value = 3124.75
$ws.Range("I137").Value = 2866.3333
$ws.Range("K137").Value = 8598.999899999999
$ws.Range("M137").Value = -6048.999899999999

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 979.6
$ws.Range("I4").Value = 724.5
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 724.5
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -608.5
$ws.Range("N4").Value = -2232
# Row 61
$ws.Range("H61").Value = 3366.4443
$ws.Range("I61").Value = 3366.4443
$ws.Range("K61").Value = 3366.4443
$ws.Range("M61").Value = -3154.4443
# Row 92
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
# Row 95
$ws.Range("H95").Value = 341332.66
$ws.Range("J95").Value = 341332.66
$ws.Range("L95").Value = 341332.66
$ws.Range("N95").Value = -346824.66
# Row 132
$ws.Range("H132").Value = 2071.2307
$ws.Range("I132").Value = 1720.5454
$ws.Range("K132").Value = 5161.6362
$ws.Range("M132").Value = -2631.6362
# Row 136
$ws.Range("H136").Value = 3366.4443
$ws.Range("I136").Value = 3366.4443
$ws.Range("K136").Value = 10099.3329
$ws.Range("M136").Value = -7549.332900000001

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1503.5555
$ws.Range("I107").Value = 1232
$ws.Range("K107").Value = 1232
$ws.Range("M107").Value = 688
# Row 134
$ws.Range("H134").Value = 3426.2144
$ws.Range("I134").Value = 2330.5833
$ws.Range("K134").Value = 6991.749899999999
$ws.Range("M134").Value = -4456.749899999999

$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 13999.5
$ws.Range("J28").Value = 13999.5
$ws.Range("L28").Value = 13999.5
$ws.Range("N28").Value = -14489.5
# Row 31
$ws.Range("H31").Value = 4267.857
$ws.Range("I31").Value = 1860.125
$ws.Range("K31").Value = 1860.125
$ws.Range("M31").Value = -1565.125
# Row 34
$ws.Range("H34").Value = 4267.857
$ws.Range("I34").Value = 1860.125
$ws.Range("K34").Value = 1860.125
$ws.Range("M34").Value = -1658.125
# Row 51
$ws.Range("H51").Value = 25000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26472
# Row 61
$ws.Range("H61").Value = 25000
$ws.Range("J61").Value = 25000
$ws.Range("L61").Value = 25000
$ws.Range("N61").Value = -25696

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 714661.3
$ws.Range("I4").Value = 1000501.2
$ws.Range("J4").Value = 61.5
$ws.Range("K4").Value = 3001503.6
$ws.Range("L4").Value = 184.5
$ws.Range("M4").Value = -3001391.6
$ws.Range("N4").Value = -408.5
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 19122.375
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 19122.375
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -20180.375
# Row 41
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
# Row 58
$ws.Range("H58").Value = 19999.5
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
# Row 102
$ws.Range("H102").Value = 2119.2727
$ws.Range("I102").Value = 2153.3
$ws.Range("J102").Value = 1779
$ws.Range("K102").Value = 2153.3
$ws.Range("L102").Value = 1779
$ws.Range("M102").Value = -531.3000000000002
$ws.Range("N102").Value = -5023
# Row 132
$ws.Range("H132").Value = 4599.467
$ws.Range("I132").Value = 4559.5
$ws.Range("J132").Value = 4679.4
$ws.Range("K132").Value = 13678.5
$ws.Range("L132").Value = 14038.2
$ws.Range("M132").Value = -11148.5
$ws.Range("N132").Value = -19098.2

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7488.1113
$ws.Range("J7").Value = 7066
$ws.Range("L7").Value = 7066
$ws.Range("N7").Value = -7290
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 68
$ws.Range("H68").Value = 3499.7778
$ws.Range("I68").Value = 2999.75
$ws.Range("J68").Value = 7500
$ws.Range("K68").Value = 2999.75
$ws.Range("L68").Value = 7500
$ws.Range("M68").Value = -2250.75
$ws.Range("N68").Value = -8998
# Row 71
$ws.Range("H71").Value = 3499.7778
$ws.Range("I71").Value = 2999.75
$ws.Range("J71").Value = 7500
$ws.Range("K71").Value = 14998.75
$ws.Range("L71").Value = 37500
$ws.Range("M71").Value = -11254.75
$ws.Range("N71").Value = -44988
# Row 126
$ws.Range("H126").Value = 7488.1113
$ws.Range("J126").Value = 7066
$ws.Range("L126").Value = 21198
$ws.Range("N126").Value = -26138
# Row 132
$ws.Range("H132").Value = 5614.2856
$ws.Range("I132").Value = 5614.2856
$ws.Range("K132").Value = 16842.8568
$ws.Range("M132").Value = -14312.8568

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 900.5
$ws.Range("J113").Value = 799
$ws.Range("L113").Value = 2397
$ws.Range("N113").Value = -6737
# Row 136
$ws.Range("H136").Value = 3558.818
$ws.Range("I136").Value = 3558.818
$ws.Range("K136").Value = 10676.454
$ws.Range("M136").Value = -8126.454000000002

Write-Host "Applied 184 cell updates and 7 clears"